$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add OC_TPR values in column B (header already present as "OC_TPR")
$ws.Range("B2").Value = 0.090909090909090898
$ws.Range("B3").Value = 0.19008264462809901
$ws.Range("B4").Value = 0.27272727272727199
$ws.Range("B5").Value = 0.330578512396694

# Adjust column widths to fit new content (column C is left untouched,
# it already has the right width from the original file)
$ws.Columns.Item(1).ColumnWidth = 12.333333333333334
$ws.Columns.Item(2).ColumnWidth = 13.666666666666666
$ws.Columns.Item(4).ColumnWidth = 14.5

# Update selection to B5 to match final cursor position
$ws.Range("B5").Select()

$wb.Save()
